$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = New-Object 'object[,]' 1,20
$row2[0,0] = "FAPs"
$row2[0,1] = "Efna5"
$row2[0,2] = "Ephb6"
$row2[0,3] = "ECs"
$row2[0,4] = 3
$row2[0,5] = 1
$row2[0,6] = 1.666083666666667
$row2[0,7] = 4.998251
$row2[0,8] = 0.6125276070882968
$row2[0,9] = 0.6125276070882968
$row2[0,10] = 2
$row2[0,11] = 0.6666666666666666
$row2[0,12] = 0.095455
$row2[0,13] = 0.286365
$row2[0,14] = 0.03201851307709132
$row2[0,15] = 0.03201851307709131
$row2[0,16] = 0.1590360164016666
$row2[0,17] = 1.431324147615
$row2[0,18] = 0.01961222319763608
$row2[0,19] = 0.01961222319763608
$ws.Range("A2:T2").Value2 = $row2

$row3 = New-Object 'object[,]' 1,20
$row3[0,0] = "FAPs"
$row3[0,1] = "Efna5"
$row3[0,2] = "Ephb6"
$row3[0,3] = "FAPs"
$row3[0,4] = 3
$row3[0,5] = 1
$row3[0,6] = 1.666083666666667
$row3[0,7] = 4.998251
$row3[0,8] = 0.6125276070882968
$row3[0,9] = 0.6125276070882968
$row3[0,10] = 3
$row3[0,11] = 1
$row3[0,12] = 1.220310333333333
$row3[0,13] = 3.660931
$row3[0,14] = 0.4093292375039861
$row3[0,15] = 0.409329237503986
$row3[0,16] = 2.033139114631222
$row3[0,17] = 18.298252031681
$row3[0,18] = 0.2507254583595937
$row3[0,19] = 0.2507254583595937
$ws.Range("A3:T3").Value2 = $row3

$row4 = New-Object 'object[,]' 1,20
$row4[0,0] = "FAPs"
$row4[0,1] = "Efna5"
$row4[0,2] = "Ephb6"
$row4[0,3] = "sCs"
$row4[0,4] = 3
$row4[0,5] = 1
$row4[0,6] = 1.666083666666667
$row4[0,7] = 4.998251
$row4[0,8] = 0.6125276070882968
$row4[0,9] = 0.6125276070882968
$row4[0,10] = 3
$row4[0,11] = 1
$row4[0,12] = 1.665478666666666
$row4[0,13] = 4.996435999999999
$row4[0,14] = 0.5586522494189227
$row4[0,15] = 0.5586522494189227
$row4[0,16] = 2.77482680371511
$row4[0,17] = 24.973441233436
$row4[0,18] = 0.342189925531067
$row4[0,19] = 0.342189925531067
$ws.Range("A4:T4").Value2 = $row4

$row5 = New-Object 'object[,]' 1,20
$row5[0,0] = "sCs"
$row5[0,1] = "Efna5"
$row5[0,2] = "Ephb6"
$row5[0,3] = "ECs"
$row5[0,4] = 3
$row5[0,5] = 1
$row5[0,6] = 1.053930333333333
$row5[0,7] = 3.161791
$row5[0,8] = 0.3874723929117032
$row5[0,9] = 0.3874723929117031
$row5[0,10] = 2
$row5[0,11] = 0.6666666666666666
$row5[0,12] = 0.095455
$row5[0,13] = 0.286365
$row5[0,14] = 0.03201851307709132
$row5[0,15] = 0.03201851307709131
$row5[0,16] = 0.1006029199683333
$row5[0,17] = 0.9054262797149999
$row5[0,18] = 0.01240628987945524
$row5[0,19] = 0.01240628987945523
$ws.Range("A5:T5").Value2 = $row5

$row6 = New-Object 'object[,]' 1,20
$row6[0,0] = "sCs"
$row6[0,1] = "Efna5"
$row6[0,2] = "Ephb6"
$row6[0,3] = "FAPs"
$row6[0,4] = 3
$row6[0,5] = 1
$row6[0,6] = 1.053930333333333
$row6[0,7] = 3.161791
$row6[0,8] = 0.3874723929117032
$row6[0,9] = 0.3874723929117031
$row6[0,10] = 3
$row6[0,11] = 1
$row6[0,12] = 1.220310333333333
$row6[0,13] = 3.660931
$row6[0,14] = 0.4093292375039861
$row6[0,15] = 0.409329237503986
$row6[0,16] = 1.286122076380111
$row6[0,17] = 11.575098687421
$row6[0,18] = 0.1586037791443924
$row6[0,19] = 0.1586037791443923
$ws.Range("A6:T6").Value2 = $row6

$row7 = New-Object 'object[,]' 1,20
$row7[0,0] = "sCs"
$row7[0,1] = "Efna5"
$row7[0,2] = "Ephb6"
$row7[0,3] = "sCs"
$row7[0,4] = 3
$row7[0,5] = 1
$row7[0,6] = 1.053930333333333
$row7[0,7] = 3.161791
$row7[0,8] = 0.3874723929117032
$row7[0,9] = 0.3874723929117031
$row7[0,10] = 3
$row7[0,11] = 1
$row7[0,12] = 1.665478666666666
$row7[0,13] = 4.996435999999999
$row7[0,14] = 0.5586522494189227
$row7[0,15] = 0.5586522494189227
$row7[0,16] = 1.755298486319555
$row7[0,17] = 15.797686376876
$row7[0,18] = 0.2164623238878556
$row7[0,19] = 0.2164623238878556
$ws.Range("A7:T7").Value2 = $row7
